# Remove the existing _GoBack bookmark (it will be re-added later, after the
# "// ..." comment line inside the new tsconfig.json code block, per the target diff).
$d = $word.ActiveDocument

try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
    # no existing _GoBack bookmark - nothing to remove
}

# Locate the "...continuing..." paragraph (last bullet under the "Angular" heading)
# and replace it (and add many new paragraphs after it) using InsertXML so that
# paragraph styles (ListParagraph / Code), numbering (numId 8 @ ilvl 0/1), run
# formatting (colors, eastAsia font), proofErr spans, tabs, lastRenderedPageBreak
# and the relocated _GoBack bookmark are all reproduced exactly.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "…continuing…`r") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the placeholder paragraph to replace"
}

$r = $target.Range

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="8"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Add environment settings</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="8"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Ensure RUN scripts for all environment</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="8"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Make sure enough memory is allocated for the app</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="8"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>Set up budgets</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="8"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>(optionally) Disable: AOT and OPTIMIZATIONS</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="8"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Set up important </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>TypeScript</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> constraints</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>tsconfig.json</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>compileOnSave</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>false</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>compilerOptions</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:tab/>
        <w:t>// …</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>  },</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>angularCompilerOptions</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: {</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"enableI18nLegacyMessageIdFormat"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>false</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>strictInjectionParameters</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>strictInputAccessModifiers</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>strictTemplates</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>enableIvy</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>forceConsistentCasingInFileNames</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>strict</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>alwaysStrict</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>noUnusedLocals</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>false</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>noImplicitReturns</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>noImplicitAny</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>noImplicitThis</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>    </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>strictNullChecks</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
        </w:rPr>
        <w:t>"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="569CD6"/>
        </w:rPr>
        <w:t>true</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>  }</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Code"/>
        <w:ind w:left="360"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman"/>
          <w:color w:val="D4D4D4"/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)

Write-Host "Edit applied"
